$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Phases")
$ws.Range("A1").Value = "TEST"
